$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.118269205093384
$ws.Range("B1").Value = 1.91999614238739
$ws.Range("C1").Value = 4.457293033599854
$ws.Range("D1").Value = 0.297851949930191
$ws.Range("E1").Value = 0.3427042067050934
